# Update "Project 1 - Final presentation.pptx" - slide 2 ("Dataset - 1"):
# add two new bullet paragraphs ("Anonymous surveyes among students" and
# "~28'000 rows") right after the "Study on stress and depression level..."
# bullet in the content placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Find the content placeholder shape ("Segnaposto contenuto 2") on the slide
# rather than assuming a fixed shape index.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -like "*Study on stress and depression level*") {
            $targetShape = $shp
        }
    }
}

$tf = $targetShape.TextFrame
$tr = $tf.TextRange

# Locate the paragraph that ends with "...other factors." so the two new
# bullets are inserted right after it, inheriting its exact paragraph /
# run formatting (bullet glyph, font, size, ...).
$count = $tr.Paragraphs().Count
$anchorPara = $null
for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "Study on stress and depression level*") {
        $anchorPara = $para
    }
}

$anchorPara.InsertAfter("`rAnonymous surveyes among students`r~28’000 rows")
